# Trade #97 closed at 2026-02-18 00:33:50 - unknown UNKNOWN +0.000%
#
# This script updates the live trading workbook to reflect:
#   1) Trade #126 (MarketMaking) being closed via an "early_exit"
#      (Exit Price 0.98, P&L% 2.0833, P&L$ 0.02, Capital After 99.54,
#       Duration 0.16 min) on both the "All Trades" sheet and the
#       strategy-specific "MarketMaking" sheet.
#   2) A brand new OPEN trade #155 for the "momentum" strategy being
#      appended to both the "All Trades" sheet and the strategy-specific
#      "momentum" sheet.
#   3) The roll-up numbers on "Summary" and "Strategy Status" sheets being
#      refreshed to account for the newly closed trade.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    # Forces a literal text value into a cell even when the text looks like
    # a date/time (e.g. "2026-02-18", "00:33:44"), which Excel would
    # otherwise silently coerce into a date serial number. The leading
    # apostrophe forces "treat as text" semantics, and resetting the style
    # back to Normal afterwards strips the date number-format that Excel
    # still auto-applies to the cell.
    param($range, [string]$text)
    $range.Formula = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) Summary sheet roll-up metrics
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.46   # Current Capital
$summary.Range("B4").Value = 0.57      # Total P&L $
$summary.Range("B6").Value = 125       # Total Trades
$summary.Range("B7").Value = 59        # Winning Trades
$summary.Range("B9").Value = 47.2      # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.54      # Capital
$status.Range("D6").Value = 45         # Trades
$status.Range("E6").Value = -0.27      # P&L $
$status.Range("F6").Value = -0.46      # P&L %
$status.Range("G6").Value = 46.67      # Win Rate %

# ---------------------------------------------------------------------
# 3) All Trades sheet - close trade #126 (row 127) and append trade #155
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close existing MarketMaking trade (row 127 = Trade # 126)
$allTrades.Range("G127").Value = 0.98        # Exit Price
$allTrades.Range("H127").Value = "CLOSED"    # Status
$allTrades.Range("I127").Value = 2.0833      # P&L %
$allTrades.Range("J127").Value = 0.02        # P&L $
$allTrades.Range("K127").Value = 99.54       # Capital After
Set-TextCell $allTrades.Range("L127") "early_exit"   # Exit Reason
$allTrades.Range("M127").Value = 0.16        # Duration (min)

# Append brand-new OPEN momentum trade (row 156 = Trade # 155)
$allTrades.Range("A156").Value = 155
Set-TextCell $allTrades.Range("B156") "2026-02-18"
Set-TextCell $allTrades.Range("C156") "00:33:44"
Set-TextCell $allTrades.Range("D156") "momentum"
Set-TextCell $allTrades.Range("E156") "UP"
$allTrades.Range("F156").Value = 0.96
# G156 (Exit Price) intentionally left blank - trade is still OPEN
Set-TextCell $allTrades.Range("H156") "OPEN"
$allTrades.Range("I156").Value = 0
$allTrades.Range("J156").Value = 0
$allTrades.Range("K156").Value = 99.23374292899115
# L156 (Exit Reason) intentionally left blank - trade is still OPEN
$allTrades.Range("M156").Value = 0
$allTrades.Range("N156").Value = 0
$allTrades.Range("O156").Value = 0
$allTrades.Range("P156").Value = 0.9
Set-TextCell $allTrades.Range("Q156") "Upward momentum: 1.980% over 10 samples"

# ---------------------------------------------------------------------
# 4) momentum sheet - append trade #155 (row 39)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Range("A39").Value = 155
Set-TextCell $momentum.Range("B39") "2026-02-18"
Set-TextCell $momentum.Range("C39") "00:33:44"
Set-TextCell $momentum.Range("D39") "momentum"
Set-TextCell $momentum.Range("E39") "UP"
$momentum.Range("F39").Value = 0.96
# G39 (Exit Price) intentionally left blank - trade is still OPEN
Set-TextCell $momentum.Range("H39") "OPEN"
$momentum.Range("I39").Value = 0
$momentum.Range("J39").Value = 0
$momentum.Range("K39").Value = 99.23374292899115
$momentum.Range("L39").Value = 0
$momentum.Range("M39").Value = 0
$momentum.Range("N39").Value = 0.9
Set-TextCell $momentum.Range("O39") "Upward momentum: 1.980% over 10 samples"
# P39 (Exit Reason) intentionally left blank - trade is still OPEN
$momentum.Range("Q39").Value = 0

# ---------------------------------------------------------------------
# 5) MarketMaking sheet - close trade #126 (row 47)
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

$marketMaking.Range("G47").Value = 0.98        # Exit Price
$marketMaking.Range("H47").Value = "CLOSED"    # Status
$marketMaking.Range("I47").Value = 2.0833      # P&L %
$marketMaking.Range("J47").Value = 0.02        # P&L $
$marketMaking.Range("K47").Value = 99.54       # Capital After
# L47 (Entry Slippage) and M47 (Exit Slippage) unchanged
Set-TextCell $marketMaking.Range("P47") "early_exit"   # Exit Reason
$marketMaking.Range("Q47").Value = 0.16        # Duration (min)

Write-Host "Applied live_trading_results update for trade #126 close + trade #155 open."
